$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

# Update existing visit counts
$ws.Range("C2").Value = 117
$ws.Range("C3").Value = 20

# New card-verification rows appended to the "Sheet" worksheet
$ws.Range("A5").Value = "01FE5197502AC472"
$ws.Range("B5").Value = "S4 phone"
$ws.Range("C5").Value = 1

$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "00000000"
$ws.Range("A6").ClearFormats()
$ws.Range("B6").Value = "s4.2"
$ws.Range("C6").Value = 1

$ws.Range("A7").Value = "08CE7849"
$ws.Range("B7").Value = "s4.2"
$ws.Range("C7").Value = 1

$ws.Range("A8").Value = "08A088E0"
$ws.Range("B8").Value = "hello there"
$ws.Range("C8").Value = 1

$ws.Range("A9").Value = "08AB5506"
$ws.Range("B9").Value = "lklk"
$ws.Range("C9").Value = 1

$ws.Range("A10").Value = "089BFDF0"
$ws.Range("B10").Value = "gary"
$ws.Range("C10").Value = 1

$ws.Range("A11").Value = "082881DB"
$ws.Range("B11").Value = "nihao"
$ws.Range("C11").Value = 1

$ws.Range("A12").Value = "01FE940C1FD75197"
$ws.Range("B12").Value = "ni"
$ws.Range("C12").Value = 1

# Update the monthly stat sheet's total formula to cover the new rows
$ws2 = $wb.Worksheets.Item("Monthly_STAT")
$ws2.Range("C3").Formula = "=SUM(sheet!E2:E12)"
